# Cover Letter edits: re-target the letter from CIBC to BASL Inc.,
# tweak the skills paragraph, drop the Peel/Encore paragraph, and
# trim the Dean's List paragraph.

$d = $word.ActiveDocument

# wdReplaceOne = 1 ; wdFindContinue (wrap) = 1
$wdReplaceOne = 1

# --- Paragraph 2: sender's employer address block -------------------
$p = $d.Paragraphs.Item(2).Range
$p.Find.Execute("CIBC", $true, $false, $false, $false, $false, $true, 1, $false, "BASL Inc.", $wdReplaceOne) | Out-Null

$p = $d.Paragraphs.Item(2).Range
$p.Find.Execute("81 Bay Street", $true, $false, $false, $false, $false, $true, 1, $false, "662 Blue Forest Hill", $wdReplaceOne) | Out-Null

$p = $d.Paragraphs.Item(2).Range
$p.Find.Execute("Toronto, ON M5J 1E6", $true, $false, $false, $false, $false, $true, 1, $false, "Burlington, ON L7L 4H3", $wdReplaceOne) | Out-Null

# --- Paragraph 3: salutation -----------------------------------------
$p = $d.Paragraphs.Item(3).Range
$p.Find.Execute("To the CIBC Technology Team", $true, $false, $false, $false, $false, $true, 1, $false, "To Whom It May Concern,", $wdReplaceOne) | Out-Null

# --- Paragraph 4: opening paragraph / position applied for -----------
$p = $d.Paragraphs.Item(4).Range
$p.Find.Execute("Application/Software Developer Co-op position.", $true, $false, $false, $false, $false, $true, 1, $false, "Full Stack Designer Co-op position at BASL Inc.", $wdReplaceOne) | Out-Null

# --- Paragraph 5: languages / skills paragraph ------------------------
$p = $d.Paragraphs.Item(5).Range
$p.Find.Execute("JavaScript, PHP, and HTML/CSS", $true, $false, $false, $false, $false, $true, 1, $false, "JavaScript, PHP, C#, and HTML/CSS", $wdReplaceOne) | Out-Null

$p = $d.Paragraphs.Item(5).Range
$p.Find.Execute("using JavaFX and have embraced", $true, $false, $false, $false, $false, $true, 1, $false, "using JavaFX and C# .NET and have embraced", $wdReplaceOne) | Out-Null

# --- Paragraphs 7+8: drop the Peel District/Encore Audio Visual paragraph
#     and fold the trimmed Dean's List paragraph into its place. -------
$apostrophe = [char]0x2019
$p7 = $d.Paragraphs.Item(7).Range
$p8 = $d.Paragraphs.Item(8).Range
$combined = $d.Range($p7.Start, $p8.End - 1)
$combined.Text = "I am also a proud member of the Dean" + $apostrophe + "s List at Mohawk College for both my first and second semesters, which reflects my passion for technology. I am excited about the opportunity to leverage my skills and experiences for your Technology Team and contribute to real-world projects."
# the old paragraph 8 mark is now an empty trailing paragraph - remove it
$d.Paragraphs.Item(8).Range.Delete() | Out-Null

# "Technology Team " -> "team " (keep the lowercase "t"+"eam " split implied by the source)
$p = $d.Paragraphs.Item(7).Range
$p.Find.Execute("your Technology Team and", $true, $false, $false, $false, $false, $true, 1, $false, "your team and", $wdReplaceOne) | Out-Null

Write-Output $d.Content.Text
